$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 89706
$ws.Range("B3").Value = 90857
$ws.Range("A4").Value = 112393218
$ws.Range("B4").Value = 89517
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = "Vedticka"
$ws.Range("G4").Value = "Fuscoporia viticola"
$ws.Range("H4").Value = "(Schwein.) Murrill"
$ws.Range("Q4").Value = 509697
$ws.Range("R4").Value = 6814179
$ws.Range("Z4").Value = "14:23"
$ws.Range("AB4").Value = "14:23"
$ws.Range("A5").Value = 112393151
$ws.Range("B5").Value = 88637
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 1962
$ws.Range("F5").Value = "Vaddporing"
$ws.Range("G5").Value = "Anomoporia kamtschatica"
$ws.Range("H5").Value = "(Parmasto) Bondartseva"
$ws.Range("Q5").Value = 509729
$ws.Range("R5").Value = 6814194
$ws.Range("Z5").Value = "14:19"
$ws.Range("AB5").Value = "14:19"
$ws.Range("A6").Value = 112344694
$ws.Range("B6").Value = 90795
$ws.Range("E6").Value = 6055
$ws.Range("F6").Value = "Spadskinn"
$ws.Range("G6").Value = "Stereopsis vitellina"
$ws.Range("H6").Value = "(S.Lundell) D.A.Reid"
$ws.Range("P6").Value = "Kanaltjärnen (Kanaltjärnen), Dlr"
$ws.Range("Q6").Value = 510393
$ws.Range("R6").Value = 6813663
$ws.Range("Z6").Value = "12:52"
$ws.Range("AB6").Value = "12:52"
$ws.Range("A7").Value = 112343580
$ws.Range("B7").Value = 89706
$ws.Range("E7").Value = 1503
$ws.Range("F7").Value = "Gräddporing"
$ws.Range("G7").Value = "Sidera lenis"
$ws.Range("H7").Value = "(P.Karst.) Miettinen"
$ws.Range("P7").Value = "Snottabo (Snottabo), Dlr"
$ws.Range("Q7").Value = 510660
$ws.Range("R7").Value = 6813980
$ws.Range("Z7").Value = "11:37"
$ws.Range("AB7").Value = "11:37"
$ws.Range("A8").Value = 112343568
$ws.Range("B8").Value = 90860
$ws.Range("E8").Value = 232140
$ws.Range("F8").Value = "Tajgataggsvamp"
$ws.Range("G8").Value = "Phellodon secretus"
$ws.Range("H8").Value = "Niemelä & Kinnunen"
$ws.Range("A9").Value = 112343724
$ws.Range("B9").Value = 89820
$ws.Range("D9").Value = "EN"
$ws.Range("E9").Value = 71
$ws.Range("F9").Value = "Urskogsporing"
$ws.Range("G9").Value = "Neoantrodia infirma"
$ws.Range("H9").Value = "(Renvall & Niemelä) Audet"
$ws.Range("Q9").Value = 510716
$ws.Range("R9").Value = 6814045
$ws.Range("Z9").Value = "11:44"
$ws.Range("AB9").Value = "11:44"
$ws.Range("A10").Value = 112343573
$ws.Range("B10").Value = 90795
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 6055
$ws.Range("F10").Value = "Spadskinn"
$ws.Range("G10").Value = "Stereopsis vitellina"
$ws.Range("H10").Value = "(S.Lundell) D.A.Reid"
$ws.Range("Q10").Value = 510660
$ws.Range("R10").Value = 6813980
$ws.Range("Z10").Value = "11:37"
$ws.Range("AB10").Value = "11:37"
$ws.Range("B11").Value = 90795
